$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-8 (the "CHIMANIMANI" institution block) were re-sorted into a new
# order, carrying each institution's Masculino/Feminino/TOTAL counts with it.
# Row 9 (GORONGOSA) is unchanged.

$ws.Range("B2").Value = "PARQUE DE CHIMANIMANI"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

$ws.Range("B3").Value = "UCM"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1

$ws.Range("B4").Value = "MICAIA"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 6

$ws.Range("B5").Value = "ISPM"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("B6").Value = "ITAM"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1

$ws.Range("B7").Value = "UNIZAMBEZE"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

$ws.Range("B8").Value = "SDAE SUSSUNDENGA"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1
